# Week 15 simulations update for the Texans "Players Data" workbook.
# Sheet 1 = Rushing, Sheet 2 = Receiving.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Rushing
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Updated stats for existing players
$rushing.Range("D3").Value = 2
$rushing.Range("E3").Value = 1
$rushing.Range("F3").Value = 1

$rushing.Range("C5").Value = 23
$rushing.Range("D5").Value = 22
$rushing.Range("E5").Value = 12
$rushing.Range("F5").Value = 7

$rushing.Range("C7").Value = 10
$rushing.Range("D7").Value = 4

# Row 10 used to be J.Akins - now it's the newly added P.Dorsett
$rushing.Range("B10").Value = "P.Dorsett"
$rushing.Range("C10").Value = 1
$rushing.Range("D10").Value = 0
$rushing.Range("E10").Value = 0
$rushing.Range("F10").Value = 0

# New row 11: J.Akins re-added at the bottom of the list
$rushing.Range("A11").Value = 9
$rushing.Range("B11").Value = "J.Akins"
$rushing.Range("C11").Value = 0
$rushing.Range("D11").Value = 0
$rushing.Range("E11").Value = 1
$rushing.Range("F11").Value = 0

# Match the bold/bordered/centered style used by the rest of column A
$rushing.Range("A2").Copy()
$rushing.Range("A11").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet 2: Receiving
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("C3").Value = 17
$receiving.Range("D3").Value = 12

$receiving.Range("C5").Value = 9
$receiving.Range("D5").Value = 7

$receiving.Range("C7").Value = 87
$receiving.Range("D7").Value = 67
$receiving.Range("E7").Value = 24
$receiving.Range("G7").Value = 9
$receiving.Range("H7").Value = 6

$receiving.Range("C8").Value = 20
$receiving.Range("D8").Value = 12

$receiving.Range("C12").Value = 29
$receiving.Range("D12").Value = 17
$receiving.Range("E12").Value = 10
$receiving.Range("F12").Value = 5
$receiving.Range("G12").Value = 4

# Row 14 used to be J.Akins - now it's the newly added D.Davis
$receiving.Range("B14").Value = "D.Davis"
$receiving.Range("C14").Value = 1
$receiving.Range("D14").Value = 1
$receiving.Range("E14").Value = 0
$receiving.Range("F14").Value = 0
$receiving.Range("G14").Value = 0
$receiving.Range("H14").Value = 0

# Row 15 used to be P.Brown - now it's the newly added P.Dorsett
$receiving.Range("B15").Value = "P.Dorsett"
$receiving.Range("C15").Value = 2
$receiving.Range("D15").Value = 1
$receiving.Range("E15").Value = 1
$receiving.Range("F15").Value = 0

# Row 16 used to be B.Jordan - now it's J.Akins (shifted down)
$receiving.Range("B16").Value = "J.Akins"
$receiving.Range("C16").Value = 28
$receiving.Range("D16").Value = 21
$receiving.Range("E16").Value = 1
$receiving.Range("G16").Value = 4
$receiving.Range("H16").Value = 1

# Row 17 used to be A.Auclair - now it's P.Brown (shifted down)
$receiving.Range("B17").Value = "P.Brown"
$receiving.Range("C17").Value = 24
$receiving.Range("D17").Value = 16
$receiving.Range("E17").Value = 3
$receiving.Range("F17").Value = 2

# New row 18: B.Jordan re-added
$receiving.Range("A18").Value = 16
$receiving.Range("B18").Value = "B.Jordan"
$receiving.Range("C18").Value = 20
$receiving.Range("D18").Value = 14
$receiving.Range("E18").Value = 3
$receiving.Range("F18").Value = 1
$receiving.Range("G18").Value = 4
$receiving.Range("H18").Value = 4

# New row 19: A.Auclair re-added
$receiving.Range("A19").Value = 17
$receiving.Range("B19").Value = "A.Auclair"
$receiving.Range("C19").Value = 3
$receiving.Range("D19").Value = 3
$receiving.Range("E19").Value = 0
$receiving.Range("F19").Value = 0
$receiving.Range("G19").Value = 1
$receiving.Range("H19").Value = 1

# Match the bold/bordered/centered style used by the rest of column A
$receiving.Range("A2").Copy()
$receiving.Range("A18").PasteSpecial(-4122)
$receiving.Range("A2").Copy()
$receiving.Range("A19").PasteSpecial(-4122)
